$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'302.31"
$ws.Range("E2").Value = "'-0.64%"

$ws.Range("D3").Value = "'36.95"
$ws.Range("E3").Value = "'3.92%"

$ws.Range("D4").Value = "'4.976"
$ws.Range("E4").Value = "'-1.90%"

$ws.Range("D5").Value = "'0.07713"
$ws.Range("E5").Value = "'-1.26%"

$ws.Range("D6").Value = "'2.086"
$ws.Range("E6").Value = "'-7.88%"

$ws.Range("D7").Value = "'7.936"
$ws.Range("E7").Value = "'-1.87%"

$ws.Range("B8").Value = "MXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D8").Value = "'0.9144"
$ws.Range("E8").Value = "'-1.59%"

$ws.Range("B9").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C9").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D9").Value = "'0.09604"
$ws.Range("E9").Value = "'2.32%"

$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "'0.1842"
$ws.Range("E10").Value = "'0.63%"

$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "'0.08549"
$ws.Range("E11").Value = "'-0.07%"

$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").Value = "'0.03508"
$ws.Range("E12").Value = "'-2.84%"

$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D13").Value = "'0.09962"
$ws.Range("E13").Value = "'-0.09%"

$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D14").Value = "'0.001471"
$ws.Range("E14").Value = "'-0.57%"

$ws.Range("B15").Value = "TigerCash"
$ws.Range("C15").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D15").Value = "'0.005700"
$ws.Range("E15").Value = "'-0.30%"

$ws.Range("B16").Value = "LEO"
$ws.Range("C16").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D16").Value = "'3.468"
$ws.Range("E16").Value = "'-0.32%"

$ws.Range("B17").Value = "GateToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D17").Value = "'4.013"
$ws.Range("E17").Value = "'-0.83%"

$ws.Range("D18").Value = "'2.220"

$ws.Range("D19").Value = "'0.3383"
$ws.Range("E19").Value = "'-0.69%"

$ws.Range("D20").Value = "'0.1326"
$ws.Range("E20").Value = "'0.30%"

$ws.Range("D21").Value = "'4.759"
$ws.Range("E21").Value = "'4.81%"

$ws.Range("D22").Value = "'0.2200"
$ws.Range("E22").Value = "'-1.64%"

$ws.Range("D23").Value = "'0.04589"
$ws.Range("E23").Value = "'-1.66%"

$ws.Range("D24").Value = "'0.005104"
$ws.Range("E24").Value = "'12.36%"

$ws.Range("D25").Value = "'0.001231"
$ws.Range("E25").Value = "'-0.07%"

$ws.Range("D26").Value = "'0.0001403"
$ws.Range("E26").Value = "'7.86%"

$ws.Range("D39").Value = "'0.01761"
$ws.Range("E39").Value = "'-0.88%"

$ws.Range("D40").Value = "'0.04602"
$ws.Range("E40").Value = "'-2.36%"

$ws.Range("D41").Value = "'0.007459"
$ws.Range("E41").Value = "'-5.81%"

$ws.Range("D42").Value = "'0.1389"
$ws.Range("E42").Value = "'-2.29%"

$ws.Range("D43").Value = "'0.007722"
$ws.Range("E43").Value = "'-3.44%"

$ws.Range("D44").Value = "'0.002164"
$ws.Range("E44").Value = "'-5.55%"

$ws.Range("D45").Value = "'0.01033"
$ws.Range("E45").Value = "'13.94%"

$ws.Range("D46").Value = "'0.00006322"
$ws.Range("E46").Value = "'1.97%"

$ws.Range("D47").Value = "'0.00000000750"
$ws.Range("E47").Value = "'-0.08%"

$ws.Range("D48").Value = "'0.0005801"
$ws.Range("E48").Value = "'0.00%"

$ws.Range("D49").Value = "'37.66"
$ws.Range("E49").Value = "'601.39%"

$ws.Range("D50").Value = "'0.002000"
$ws.Range("E50").Value = "'-25.71%"

$ws.Range("D51").Value = "'0.00002100"
$ws.Range("E51").Value = "'-0.08%"
